# Workbook edit: update the notes section on Sheet1
#  - Row 6 text is shortened to drop the "recharge MC" clause
#  - A new row (7) is inserted with a note about charging area not being
#    included in the algorithmic MC, merged A7:D7, matching the style used
#    by the other note rows
#  - Row 8 keeps its existing text (recharge-area dimensions)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shorten the note in row 6 (drop the "recharge MC" clause)
$ws.Range("A6").Value = "for all test instances, the last MC is rework MC "

# Add a new note row, reusing the same formatting as the row above it
$ws.Range("A6:D6").Copy()
[void]$ws.Range("A7:D7").PasteSpecial(-4122)
[void]$ws.Range("A7:D7").Merge()
$ws.Range("A7").Value = "Charging area not included in algorithmic MC"

[void]$ws.Range("D17").Select()
